# "long sub str  with dp"
# Adds a new row (row 3) to the "动态规划" (Dynamic Programming) sheet
# describing the "Longest Increasing Subsequence" problem, solved with DP.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new data row ----------------------------------------------------
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 300
$ws.Range("C3").Value = "给定一个无序的整数数组，找到其中最长上升子序列的长度。 `n示例: `n输入: [10,9,2,5,3,7,101,18]`n输出: 4 `n解释: 最长的上升子序列是 [2,3,7,101]，它的长度是 4。 `n说明: `n可能会有多种最长上升子序列的组合，你只需要输出对应的长度即可。 `n你算法的时间复杂度应该为 O(n2) 。 "
$ws.Range("D3").Value = "1 dp[i]代表以nums[i]结尾的子序列的最大长度`n2 dp[0]=1,以num[0]结尾的子序列即元素本身，最大长度是1`n3  例如：计算dp[5]，就需要计算dp[4]，即以4结尾的最大子序列长度`n4 dp中存储各个数字对应的最大递增序列长度，需要遍历数组，获取最大长度`""
$ws.Range("E3").Value = "最长子序列"
$ws.Range("F3").Value = "O(n*n)"
$ws.Range("G3").Value = "O(n)"

# D3 previously used its own (italic Monaco) style; make it match the
# rest of the row (same formatting as the other data cells).
$ws.Range("C3").Copy() | Out-Null
$ws.Range("D3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# The row now holds a lot of wrapped text, so it needs to be much taller.
$ws.Rows.Item(3).RowHeight = 242

# Move the sheet's selection like the author left it.
$ws.Range("E5").Select() | Out-Null
